$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "time_taken", styled the same as the other header cells (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in time_taken values for data rows 2-31 (plain text, unstyled like the rest of the data cells)
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:39:30.319400"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:39:30.319412"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:39:30.319415"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:39:30.319418"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:39:30.319421"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:39:30.319423"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:39:30.319426"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:39:30.319428"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:39:30.319431"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:39:30.319434"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:39:30.319437"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:39:30.319439"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:39:30.319442"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:39:30.319444"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:39:30.319447"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:39:30.319450"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:39:30.319452"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:39:30.319455"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:39:30.319458"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:39:30.319460"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:39:30.319463"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:39:30.319465"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:39:30.319468"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:39:30.319470"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:39:30.319473"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:39:30.319476"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:39:30.319478"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:39:30.319481"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:39:30.319483"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:39:30.319486"
